$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Populate the new named-range data, writing cells in the exact order
# --- that reproduces the target shared-strings table layout.

# column_range -> Sheet1!$D$5:$D$9
$ws.Range("D5").Value = "column_range_1"
$ws.Range("D6").Value = "column_range_2"
$ws.Range("D7").Value = "column_range_3"
$ws.Range("D8").Value = "column_range_4"
$ws.Range("D9").Value = "column_range_5"

# grid_range -> Sheet1!$B$13:$D$14 (column-major fill order)
$ws.Range("B13").Value = "grid_range_1_1"
$ws.Range("B14").Value = "grid_range_1_2"
$ws.Range("C13").Value = "grid_range_2_1"
$ws.Range("C14").Value = "grid_range_2_2"
$ws.Range("D13").Value = "grid_range_3_1"
$ws.Range("D14").Value = "grid_range_3_2"

# row_range -> Sheet1!$B$11:$D$11
$ws.Range("B11").Value = "row_range_1"
$ws.Range("C11").Value = "row_range_2"
$ws.Range("D11").Value = "row_range_3"

# --- Defined names (workbook scope) ---
$wb.Names.Add("column_range", "Sheet1!`$D`$5:`$D`$9")
$wb.Names.Add("row_range", "Sheet1!`$B`$11:`$D`$11")
$wb.Names.Add("grid_range", "Sheet1!`$B`$13:`$D`$14")
$wb.Names.Add("non_adjacent_range", "Sheet1!`$B`$11:`$D`$11,Sheet1!`$B`$13:`$D`$14,Sheet1!`$D`$5:`$D`$9")

# --- Selection moves to H10 on Sheet1 ---
$ws.Range("H10").Select() | Out-Null

# --- Rename the default cell style from "Standard" to "Normal" ---
# (Style.Name is read-only via COM, so recreate it instead.)
$wb.Styles.Item("Standard").Delete()
$wb.Styles.Add("Normal") | Out-Null
